# The "strain" column previously held the value "KN99alpha", which was an
# error - the strain should have been "TDY451". Also remove the three empty
# columns (experimentObservations, floodmedia, inductionDelay) that had no
# data in any row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty columns (rightmost first so earlier deletions don't shift
# the column letters of the ones still to be removed).
$ws.Columns("I").Delete()   # inductionDelay (empty)
$ws.Columns("H").Delete()   # floodmedia (empty)
$ws.Columns("E").Delete()   # experimentObservations (empty)

# After the deletions, the strain column (previously F, now E) still has the
# incorrect "KN99alpha" value in every data row - correct it to "TDY451".
$ws.Range("E2:E5").Value = "TDY451"

# Reselect the columns that now contain treatment/timePoint data.
$ws.Columns("G:H").Select()
